# RELRADResultsSimpleTest.xlsx — "added DERS and fixed some things"
#
# LP4 (row 5) customer type changes from "residential" to a new "industrial"
# category, and the load-level increase for LP2/LP3/LP4/TOTAL propagates
# into their derived columns (F..M). Row 2 (LP1) is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (LP2): U [F] increases, ripple into SAIFI/SAIDI/CAIDI/EENS ---
$ws.Cells.Item(3, 6).Value  = 3.050000000000001   # F3  U
$ws.Cells.Item(3, 8).Value  = 1.386363636363636   # H3  SAIFI
$ws.Cells.Item(3, 11).Value = 640.5000000000001   # K3  SAIDI
$ws.Cells.Item(3, 12).Value = 291.1363636363637   # L3  CAIDI
$ws.Cells.Item(3, 13).Value = 1.63175             # M3  EENS

# --- Row 4 (LP3): same ripple ---
$ws.Cells.Item(4, 6).Value  = 3.8                 # F4  U
$ws.Cells.Item(4, 8).Value  = 1.727272727272727   # H4  SAIFI
$ws.Cells.Item(4, 11).Value = 798                 # K4  SAIDI
$ws.Cells.Item(4, 12).Value = 362.7272727272727   # L4  CAIDI
$ws.Cells.Item(4, 13).Value = 2.033               # M4  EENS

# --- Row 5 (LP4): customer type -> "industrial" (new shared string) ---
$ws.Cells.Item(5, 2).Value  = "industrial"        # B5  Customer type
$ws.Cells.Item(5, 6).Value  = 4.200000000000001   # F5  U
$ws.Cells.Item(5, 8).Value  = 1.909090909090909   # H5  SAIFI
$ws.Cells.Item(5, 11).Value = 4.200000000000001   # K5  SAIDI
$ws.Cells.Item(5, 12).Value = 1.909090909090909   # L5  CAIDI
$ws.Cells.Item(5, 13).Value = 2.3772              # M5  EENS

# --- Row 6 (TOTAL): recomputed aggregate ---
$ws.Cells.Item(6, 11).Value = 2.985261489698891   # K6  SAIDI total
$ws.Cells.Item(6, 12).Value = 1.356937040772223   # L6  CAIDI total
$ws.Cells.Item(6, 13).Value = 7.165450000000002   # M6  EENS total
